# Update "想去人数" (people interested) figures for the 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 6677
    $ws.Range("F5").Value = 51
    $ws.Range("F6").Value = 2080
    $ws.Range("F7").Value = 1587
    $ws.Range("F8").Value = 314
    $ws.Range("F10").Value = 467
    $ws.Range("F12").Value = 5653
}
